$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 11 ("MongoDB (normal, no compression)" block),
# shifting the existing rows (separator, MongoDB row, Oracle row, and the
# trailing footnote rows) down by one.
$ws.Rows.Item(11).Insert()

# Populate the newly inserted row 11 with the
# "Postgres BYTEA (disabled 2nd level hibernate cache, lz4 compression)" data.
$ws.Range("A11").Value = "Postgres BYTEA (disabled 2nd level hibernate cache, lz4 compression)"
$ws.Range("B11").Value = 11
$ws.Range("C11").Value = 131
$ws.Range("D11").Value = 262
$ws.Range("E11").Value = 8
$ws.Range("F11").Value = 67
$ws.Range("G11").Value = 126
$ws.Range("H11").Value = "local docker"

# Match the author's final selection position.
$ws.Range("G12").Select()
